$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove one of the two consecutive empty paragraphs that precede the
#    "Dualità" heading (collapsing them into a single empty paragraph).
# ---------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("se massimo va letta da duale a primale)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchorPara = $anchor.Paragraphs(1)
$emptyPara = $anchorPara.Next()
$emptyPara.Range.Delete()

# ---------------------------------------------------------------------
# 2) " segno delle disuguaglianze rispetto al corrispondente vincolo
#    primale" -> " segno " + "della corrispondente variabile di dominio
#    primale" (split into two runs).
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("delle disuguaglianze rispetto al corrispondente vincolo primale", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
# Pre-split this run away from the preceding "del" run before rewriting
# its text, otherwise the text assignment below would merge across that
# run boundary too.
$rng.Font.Bold = 1
$rng.Font.Bold = 0
$rng.Text = "della corrispondente variabile di dominio primale"
$splitRng = $d.Range($rng.Start, $rng.End)
$splitRng.Font.Bold = 1
$splitRng.Font.Bold = 0

# ---------------------------------------------------------------------
# 3) ", riporto lo stesso segno delle disuguaglianze " + "di dominio " +
#    "del corrispondente " + "problema" + " primale" ->
#    ", riporto lo stesso segno " + "della corrispondente variabile di
#    dominio primale" (collapsed to two runs).
# ---------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute(", riporto lo stesso segno delle disuguaglianze di dominio del corrispondente problema primale", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng2.Text = ", riporto lo stesso segno della corrispondente variabile di dominio primale"
$prefix = ", riporto lo stesso segno "
$splitPoint2 = $rng2.Start + $prefix.Length
$splitRng2 = $d.Range($splitPoint2, $rng2.End)
$splitRng2.Font.Bold = 1
$splitRng2.Font.Bold = 0
